# Capturados nuevas celdas en cada tipo de archivo
# Insert a new data row at row 2 (pushing the existing rows 2-75 down to 3-76)
# and populate it with the new "Período" concept captured for this report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 2; Excel shifts rows 2..75 down to 3..76
$ws.Rows.Item(2).Insert()

# The freshly inserted row picks up the header row's style by default.
# Re-apply the standard data-row formatting (same as the rows below it)
# by copying the format from what is now row 3 (the former row 2).
$ws.Range("A3:E3").Copy()
$ws.Range("A2:E2").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows.Item(2).RowHeight = $ws.Rows.Item(3).RowHeight

# Populate the new row with the captured "Período" concept data
$ws.Range("A2").Value = "Período"
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 2019
$ws.Range("E2").Value = "INVERSIONES ORTIZ VASQUEZ HERMANOS S A S"
